$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-28 11:13:18"

$wsZhCn.Range("H2").Value = "2016-08-28 11:13:13"
$wsZhCn.Range("K2").Value = "2016-08-28 11:13:29"

$wsDeDe.Range("H2").Value = "2016-08-28 11:13:18"
$wsDeDe.Range("K2").Value = "2016-08-28 11:13:37"
